# NBA2324.xlsx — append the next batch of 9 completed games to Sheet1
# (rows 734-742), then move the view/selection the same way the author
# left it (scrolled down, G735 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Away team, B=Away Pts, C=Home team, D=Home Pts, E=Overtime,
#          F=Attend., G=Arena, H=Win, I=Loss
$newGames = @(
    @{ Row=734; Away="Orlando Magic";          AwayPts=111; Home="Detroit Pistons";         HomePts=99;  OT="No"; Attend=17832; Arena="Little Caesars Arena"; Win="Orlando Magic";          Loss="Detroit Pistons" }
    @{ Row=735; Away="Phoenix Suns";            AwayPts=140; Home="Washington Wizards";      HomePts=112; OT="No"; Attend=17832; Arena="Capital One Arena";    Win="Phoenix Suns";            Loss="Washington Wizards" }
    @{ Row=736; Away="Memphis Grizzlies";       AwayPts=91;  Home="Boston Celtics";          HomePts=131; OT="No"; Attend=17832; Arena="TD Garden";            Win="Boston Celtics";          Loss="Memphis Grizzlies" }
    @{ Row=737; Away="Indiana Pacers";          AwayPts=115; Home="Charlotte Hornets";       HomePts=99;  OT="No"; Attend=17832; Arena="Spectrum Center";      Win="Indiana Pacers";          Loss="Charlotte Hornets" }
    @{ Row=738; Away="Los Angeles Clippers";    AwayPts=103; Home="Miami Heat";              HomePts=95;  OT="No"; Attend=17832; Arena="Kaseya Center";        Win="Los Angeles Clippers";    Loss="Miami Heat" }
    @{ Row=739; Away="Houston Rockets";         AwayPts=90;  Home="Minnesota Timberwolves";  HomePts=111; OT="No"; Attend=17832; Arena="Target Center";        Win="Minnesota Timberwolves";  Loss="Houston Rockets" }
    @{ Row=740; Away="Toronto Raptors";         AwayPts=127; Home="Oklahoma City Thunder";   HomePts=135; OT="No"; Attend=17832; Arena="Paycom Center";        Win="Oklahoma City Thunder";   Loss="Toronto Raptors" }
    @{ Row=741; Away="Milwaukee Bucks";         AwayPts=108; Home="Utah Jazz";               HomePts=123; OT="No"; Attend=17832; Arena="Delta Center";         Win="Utah Jazz";               Loss="Milwaukee Bucks" }
    @{ Row=742; Away="Portland Trail Blazers";  AwayPts=103; Home="Denver Nuggets";          HomePts=112; OT="No"; Attend=17832; Arena="Ball Arena";           Win="Denver Nuggets";          Loss="Portland Trail Blazers" }
)

foreach ($g in $newGames) {
    $r = $g.Row
    $ws.Cells.Item($r, 1).Value = $g.Away
    $ws.Cells.Item($r, 2).Value = $g.AwayPts
    $ws.Cells.Item($r, 3).Value = $g.Home
    $ws.Cells.Item($r, 4).Value = $g.HomePts
    $ws.Cells.Item($r, 5).Value = $g.OT
    $ws.Cells.Item($r, 6).Value = $g.Attend
    $ws.Cells.Item($r, 7).Value = $g.Arena
    $ws.Cells.Item($r, 8).Value = $g.Win
    $ws.Cells.Item($r, 9).Value = $g.Loss

    # Points columns (B, D) carry the thousands-separator number format
    # used throughout the sheet; other columns stay General/text.
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0"
}

# Match the author's final scroll position / selection (topLeftCell A709,
# active cell G735) as closely as the host lets us.
$win = $excel.ActiveWindow
$win.ScrollRow = 709
$win.ScrollColumn = 1
$null = $ws.Range("G735").Select()
